$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.104.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.580.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.581"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.586.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0995"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.136"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.062.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.075.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.609.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000131"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "332.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.419"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0727"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "36.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.826"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.814"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("E40").Value = "  -4.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "279.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.589"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0945"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0527"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0225"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.924.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.46%  "
